$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price-table refresh (Mon Jan 16 17:24:43 UTC 2023 GitHub Actions run).
# Only the "Price" (D) and "Volume(1h)" (E) columns move; everything else is
# untouched. The sheet stores these as plain text (e.g. "297.39", "-1.11%"), so
# NumberFormat is forced to "@" (Text) before each write and the cell style is
# reset to "Normal" afterwards - this keeps the values as text instead of being
# auto-parsed into numbers/percentages, while leaving no extra cell formatting
# behind.
$updates = [ordered]@{
    "D2" = "297.27"
    "E2" = "-1.28%"
    "D3" = "31.67"
    "E3" = "0.47%"
    "D4" = "5.060"
    "E4" = "-1.91%"
    "D5" = "0.08056"
    "E5" = "8.69%"
    "D6" = "2.601"
    "E6" = "15.33%"
    "E7" = "-1.44%"
    "D8" = "3.819"
    "E8" = "1.26%"
    "D9" = "0.9231"
    "E9" = "-0.60%"
    "D10" = "0.1757"
    "E10" = "2.36%"
    "D11" = "0.07421"
    "E11" = "-2.80%"
    "D12" = "0.08875"
    "E12" = "8.53%"
    "D13" = "0.03034"
    "E13" = "0.15%"
    "E14" = "0.82%"
    "D15" = "0.001508"
    "E15" = "0.54%"
    "D16" = "0.006000"
    "E16" = "-2.67%"
    "E17" = "2.51%"
    "D18" = "2.249"
    "E18" = "0.84%"
    "D19" = "0.3258"
    "E19" = "0.33%"
    "D20" = "0.1336"
    "E20" = "-0.02%"
    "D21" = "3.994"
    "E21" = "-14.08%"
    "D22" = "0.1648"
    "E22" = "4.15%"
    "D23" = "0.04587"
    "E23" = "-1.36%"
    "D24" = "0.001242"
    "E24" = "1.56%"
    "D25" = "0.004447"
    "E25" = "-0.82%"
    "E26" = "-7.78%"
    "D27" = "0.0003407"
    "E27" = "82.04%"
    "D39" = "0.01777"
    "E39" = "2.54%"
    "D40" = "0.04493"
    "E40" = "-0.81%"
    "D41" = "0.006726"
    "E41" = "-5.69%"
    "E42" = "0.10%"
    "D43" = "0.002206"
    "E43" = "0.82%"
    "D44" = "0.009839"
    "E44" = "-6.65%"
    "D45" = "0.00006466"
    "E45" = "2.90%"
    "E46" = "-0.18%"
    "D47" = "0.008735"
    "E47" = "24.84%"
    "D49" = "0.00002096"
    "E49" = "-0.18%"
    "E50" = "-0.11%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

